# Refresh the crypto price/volume table (GitHub Actions daily pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '34.477.15'
$ws.Range("E2").Value = '  +1.07%  '

# Row 3
$ws.Range("D3").Value = '1.797.33'

# Row 4
$ws.Range("E4").Value = '  -0.11%  '

# Row 5
$ws.Range("D5").Value = "'227.18"
$ws.Range("E5").Value = '  +0.34%  '

# Row 6
$ws.Range("E6").Value = '  +1.63%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("D8").Value = "'32.48"
$ws.Range("E8").Value = '  +2.11%  '

# Row 9
$ws.Range("E9").Value = '  +1.60%  '

# Row 10
$ws.Range("E10").Value = '  +0.97%  '

# Row 11
$ws.Range("D11").Value = "'0.0949"
$ws.Range("E11").Value = '  +0.46%  '

# Row 12
$ws.Range("D12").Value = '2.057.58'
$ws.Range("E12").Value = '  +0.56%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = "'11.10"
$ws.Range("E13").Value = '  -0.37%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.792.35'
$ws.Range("E14").Value = '  +0.51%  '

# Row 15
$ws.Range("D15").Value = "'0.637"
$ws.Range("E15").Value = '  +2.83%  '

# Row 16
$ws.Range("D16").Value = '34.426.06'
$ws.Range("E16").Value = '  +1.05%  '

# Row 17
$ws.Range("E17").Value = '  +1.93%  '

# Row 18
$ws.Range("D18").Value = "'68.61"
$ws.Range("E18").Value = '  +0.93%  '

# Row 19
$ws.Range("D19").Value = "'247.14"
$ws.Range("E19").Value = '  +0.78%  '

# Row 20
$ws.Range("E20").Value = '  +3.16%  '

# Row 21
$ws.Range("D21").Value = "'11.19"
$ws.Range("E21").Value = '  +3.16%  '

# Row 22
$ws.Range("E22").Value = '  -0.10%  '

# Row 23
$ws.Range("E23").Value = '  +1.60%  '

# Row 24
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = '  +1.67%  '

# Row 25
$ws.Range("E25").Value = '  +1.00%  '

# Row 26
$ws.Range("D26").Value = "'7.27"
$ws.Range("E26").Value = '  +1.85%  '

# Row 27
$ws.Range("D27").Value = "'16.58"
$ws.Range("E27").Value = '  +1.79%  '

# Row 30
$ws.Range("E30").Value = '  +0.56%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'3.92"
$ws.Range("E31").Value = '  +8.81%  '

# Row 32
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = "'0.0523"
$ws.Range("E32").Value = '  +0.97%  '

# Row 33
$ws.Range("D33").Value = "'3.79"
$ws.Range("E33").Value = '  +3.42%  '

# Row 34
$ws.Range("E34").Value = '  +1.48%  '

# Row 35
$ws.Range("D35").Value = '1.444.90'
$ws.Range("E35").Value = '  -1.08%  '

# Row 36
$ws.Range("E36").Value = '  +7.37%  '

# Row 37
$ws.Range("E37").Value = '  +3.77%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = "'1.06"
$ws.Range("E38").Value = '  +2.23%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = "'0.0191"
$ws.Range("E39").Value = '  -0.20%  '

# Row 40
$ws.Range("E40").Value = '  +4.89%  '

# Row 41
$ws.Range("E41").Value = '  +1.38%  '

# Row 42
$ws.Range("E42").Value = '  +1.68%  '

# Row 43
$ws.Range("E43").Value = '  +2.58%  '

# Row 44
$ws.Range("D44").Value = "'13.81"
$ws.Range("E44").Value = '  +2.70%  '

# Row 45
$ws.Range("E45").Value = '  +3.25%  '

# Row 46
$ws.Range("E46").Value = '  +0.63%  '

# Row 47
$ws.Range("E47").Value = '  +0.30%  '

# Row 48
$ws.Range("D48").Value = '1.953.92'
$ws.Range("E48").Value = '  +0.34%  '

# Row 49
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0133'
$ws.Range("E49").Value = '  -1.67%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'105.75"
$ws.Range("E50").Value = '  -0.38%  '

# Row 51
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = '  -0.08%  '
